$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) counts for a few events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13231
$ws1.Range("F11").Value = 993
$ws1.Range("F27").Value = 294

# Sheet "全部类型": same underlying rows duplicated, update accordingly
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13231
$ws4.Range("F12").Value = 993
$ws4.Range("F28").Value = 294
